# Edit DataManagementSystemPrimer.docx:
#  - Append a new sentence about the validator living in the DMSTool project
#    (or a new project) to the paragraph that currently ends with
#    "...DMSValidatorControl property. "
#  - Relocate the "_GoBack" bookmark (an artifact Word leaves at the last edit
#    position) from the end of the document to the end of this newly edited
#    paragraph, matching what Word does when you stop typing there.

$d = $word.ActiveDocument

# --- Locate the insertion point robustly via Find, rather than assuming a
#     fixed paragraph index. -------------------------------------------------
$anchorText = "DMSValidatorControl property. "
$full = $d.Range(0, $d.Content.End)
$found = $full.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate anchor paragraph text to edit."
}

# Collapse to the point right after the matched text (still inside the
# paragraph, right before its paragraph mark).
$full.Collapse(0)
$insertPos = $full.Start

# --- The new sentence, split into the same run-sized chunks the final
#     document uses. ---------------------------------------------------------
$segments = @(
    " This validator can be part of the ",
    "DMSTool",
    " project, or you can create a new project, as long as the validator library files are located in the %",
    "appdata",
    "%\",
    "dmsExtensions",
    " directory for LcmsNet to use."
)

$pos = $insertPos
foreach ($seg in $segments) {
    $r = $d.Range($pos, $pos)
    $r.InsertAfter($seg)
    $pos = $pos + $seg.Length
}

# --- Move the "_GoBack" bookmark to the end of the text we just inserted. --
# Placing a bookmark exactly one position before a paragraph mark is
# unreliable, so insert a one-character placeholder first, anchor the
# bookmark right before it, then remove the placeholder again. Adding a
# bookmark named "_GoBack" automatically replaces any existing bookmark of
# that name (Word bookmark names are unique), so the old one at the end of
# the document is removed as a side effect.
$placeholderR = $d.Range($pos, $pos)
$placeholderR.InsertAfter("Z")

$bmR = $d.Range($pos, $pos)
$bmR.Bookmarks.Add("_GoBack")

$cleanupR = $d.Range($pos, $pos + 1)
$cleanupR.Text = ""
